# Applies the "IBET" tagging edit to the "Sheet3" worksheet of the
# HUMs-at-a-glance workbook.
#
# Summary of the underlying data change:
#   - Column D (and its mirrored block in column H) of the schedule table
#     used to carry "Sen Sem" / "Global" markers for the four rows at
#     A38:A41 - those markers are removed (cells cleared).
#   - The eighteen rows A45:A62 (and their H-column counterparts) gain a
#     new "IBET" marker in columns D and H, which introduces a new shared
#     string ("IBET") to the workbook.
#   - The active worksheet's selection/scroll position is updated to
#     reflect where the author was last working (cell D42).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# Clear the stale "Sen Sem" / "Global" markers in D38:D41.
$ws.Range("D38").Value = ""
$ws.Range("D39").Value = ""
$ws.Range("D40").Value = ""
$ws.Range("D41").Value = ""

# Tag rows 45-62 (columns D and H) with the new "IBET" marker.
For ($r = 45; $r -le 62; $r++) {
    $ws.Cells.Item($r, 4).Value = "IBET"
    $ws.Cells.Item($r, 8).Value = "IBET"
}

# Leave the sheet's selection where the author ended up.
$ws.Activate()
$ws.Range("D42").Select()
